$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Fri Oct 25 12:03:48 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 12:04:02 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 12:04:17 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 12:04:31 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 12:04:44 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 12:04:56 EDT 2024"
